$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column S (year 2022) to the header row (row 4), copying style from R4
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("S4").Value = 2022

# Update existing data values in row 5
$ws.Range("P5").Value = 20.5
$ws.Range("Q5").Value = 20.5
$ws.Range("R5").Value = 17.9

# Add new data value in S5, copying style from R5
$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("S5").Value = 13.5

# Update selection to match the target view state
$ws.Range("S7:S8").Select()
